$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight E2 (PREMIUM = Yes) with a yellow fill (new style/fill added by this edit)
$ws.Range("E2").Interior.Color = 65535

# Column widths
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws.Columns.Item(2).ColumnWidth = 54.166666666666664
$ws.Columns.Item(3).ColumnWidth = 65.16666666666667
$ws.Columns.Item(4).ColumnWidth = 52.166666666666664
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 59.166666666666664

# Refresh row data (rows 2-17) with the latest scraped opportunities
$ws.Range("A2").Value = "'1327954"
$ws.Range("A2").ClearFormats()
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1327954"
$ws.Range("C2").Value = "Taste Hungary | Jr. Commercial Analyst (Finance)"
$ws.Range("D2").Value = "Budapeste, Hungria"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "4 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "EATON"

$ws.Range("A3").Value = "'1327949"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1327949"
$ws.Range("C3").Value = "Taste Hungary | Football Data Analyst"
$ws.Range("D3").Value = "Budapest, Magyarország"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "3 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "ACE Advisory"

$ws.Range("A4").Value = "'1327919"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1327919"
$ws.Range("C4").Value = "Customer Representative- Intern"
$ws.Range("D4").Value = "Nugegoda, Sri Lanka"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "2 applicants"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "KAYJAY ELECTRONICS (PVT) LTD"

$ws.Range("A5").Value = "'1327871"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1327871"
$ws.Range("C5").Value = "International Relations Development Intern"
$ws.Range("D5").Value = "Jalandhar, Punjab, India"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "0 applicants"
$ws.Range("G5").Value = "3 - 6 Months"
$ws.Range("H5").Value = "Lovely Professional University"

$ws.Range("A6").Value = "'1327825"
$ws.Range("A6").ClearFormats()
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1327825"
$ws.Range("C6").Value = "Project Management Trainee"
$ws.Range("D6").Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "33 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "HILTI Panama"

$ws.Range("A7").Value = "'1327775"
$ws.Range("A7").ClearFormats()
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1327775"
$ws.Range("C7").Value = "Accelerate Romania| Programming Intern"
$ws.Range("D7").Value = "Bucharest, Romania"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "14 applicants"
$ws.Range("G7").Value = "9 - 12 Weeks"
$ws.Range("H7").Value = "AQUAsoft"

$ws.Range("A8").Value = "'1327768"
$ws.Range("A8").ClearFormats()
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1327768"
$ws.Range("C8").Value = "Accelerate Romania| Business Development Intern"
$ws.Range("D8").Value = "Bucharest, Romania"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "12 applicants"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "AQUAsoft"

$ws.Range("A9").Value = "'1326913"
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1326913"
$ws.Range("C9").Value = "Occupational Health and Safety Project Specialist (Mine Opp)"
$ws.Range("D9").Value = "Mexico City, CDMX, Mexico"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "18 applicants"
$ws.Range("G9").Value = "6 - 18 Months"
$ws.Range("H9").Value = "Sodexo Mexico"

$ws.Range("A10").Value = "'1326890"
$ws.Range("A10").ClearFormats()
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1326890"
$ws.Range("C10").Value = "Occupational Health and Safety Projects Specialist"
$ws.Range("D10").Value = "Mexico City, CDMX, Mexico"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "13 applicants"
$ws.Range("G10").Value = "6 - 18 Months"
$ws.Range("H10").Value = "Sodexo Mexico"

$ws.Range("A11").Value = "'1326448"
$ws.Range("A11").ClearFormats()
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1326448"
$ws.Range("C11").Value = "TIM Operations Assistant"
$ws.Range("D11").Value = "Panamá, Provincia de Panamá, Panamá"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "58 applicants"
$ws.Range("G11").Value = "6 - 18 Months"
$ws.Range("H11").Value = "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"

$ws.Range("A12").Value = "'1326423"
$ws.Range("A12").ClearFormats()
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1326423"
$ws.Range("C12").Value = "Sales Specialist"
$ws.Range("D12").Value = "Sincan, İstasyon, 06934 Sincan/Ankara, Türkiye"
$ws.Range("E12").Value = "No"
$ws.Range("F12").Value = "35 applicants"
$ws.Range("G12").Value = "9 - 12 Weeks"
$ws.Range("H12").Value = "Jump Start"

$ws.Range("A13").Value = "'1326270"
$ws.Range("A13").ClearFormats()
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1326270"
$ws.Range("C13").Value = "B2B Sales Manager"
$ws.Range("D13").Value = "Ciudad de Córdoba, Provincia de Córdoba, Argentina"
$ws.Range("E13").Value = "No"
$ws.Range("F13").Value = "40 applicants"
$ws.Range("G13").Value = "3 - 6 Months"
$ws.Range("H13").Value = "Bizit Global"

$ws.Range("A14").Value = "'1322493"
$ws.Range("A14").ClearFormats()
$ws.Range("B14").Value = "https://aiesec.org/opportunity/global-talent/1322493"
$ws.Range("C14").Value = "[Impact Fortaleza]- Cost & Quality Planning"
$ws.Range("D14").Value = "Castanhal, PA, Brasil"
$ws.Range("E14").Value = "No"
$ws.Range("F14").Value = "28 applicants"
$ws.Range("G14").Value = "6 - 18 Months"
$ws.Range("H14").Value = "Petruz Fruity"

$ws.Range("A15").Value = "'1321261"
$ws.Range("A15").ClearFormats()
$ws.Range("B15").Value = "https://aiesec.org/opportunity/global-talent/1321261"
$ws.Range("C15").Value = "UP Romania|Marketing General"
$ws.Range("D15").Value = "Bucharest, Romania"
$ws.Range("E15").Value = "No"
$ws.Range("F15").Value = "65 applicants"
$ws.Range("G15").Value = "9 - 12 Weeks"
$ws.Range("H15").Value = "Cluster CS"

$ws.Range("A16").Value = "'1320868"
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = "https://aiesec.org/opportunity/global-talent/1320868"
$ws.Range("C16").Value = "Accelerate Romania|Data Labeling Specialist (SERBIAN Speackers)"
$ws.Range("D16").Value = "Bucharest, Romania"
$ws.Range("E16").Value = "No"
$ws.Range("F16").Value = "7 applicants"
$ws.Range("G16").Value = "9 - 12 Weeks"
$ws.Range("H16").Value = "RepsMate"

$ws.Range("A17").Value = "'1289378"
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = "https://aiesec.org/opportunity/global-talent/1289378"
$ws.Range("C17").Value = "Medical Advisor (Spanish Speaker)"
$ws.Range("D17").Value = "İstanbul, Türkiye"
$ws.Range("E17").Value = "No"
$ws.Range("F17").Value = "116 applicants"
$ws.Range("G17").Value = "6 - 18 Months"
$ws.Range("H17").Value = "International Plus"
